$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.1443055028629487
$ws.Range("D2").Value = 0.02583992922655654
$ws.Range("E2").Value = 0.2091334449716982
$ws.Range("F2").Value = 0.6527184526782648
$ws.Range("G2").Value = 0.002417491095003315
$ws.Range("I2").Value = 0.5726874813622871
$ws.Range("L2").Value = 0.4381823451170703
$ws.Range("M2").Value = 10.88947443666467
$ws.Range("O2").Value = 2.211179229012771
$ws.Range("C3").Value = 0.1445433111169194
$ws.Range("D3").Value = 0.0249330406289161
$ws.Range("E3").Value = 0.192535924534333
$ws.Range("F3").Value = 0.6667955102320278
$ws.Range("G3").Value = 0.002422754144630197
$ws.Range("I3").Value = 0.5988095570524514
$ws.Range("L3").Value = 0.3929985317799947
$ws.Range("M3").Value = 9.560103055943898
$ws.Range("O3").Value = 2.281527690969654
$ws.Range("C4").Value = 0.1448770814070315
$ws.Range("D4").Value = 0.02439595709433462
$ws.Range("E4").Value = 0.1825285129247263
$ws.Range("F4").Value = 0.6766522063770566
$ws.Range("G4").Value = 0.002426123734616511
$ws.Range("I4").Value = 0.6160816230205661
$ws.Range("L4").Value = 0.3655212271781068
$ws.Range("M4").Value = 8.740752815026383
$ws.Range("O4").Value = 2.329103703260373
$ws.Range("C5").Value = 0.145059556483659
$ws.Range("D5").Value = 0.02418197293101088
$ws.Range("E5").Value = 0.1784949067261152
$ws.Range("F5").Value = 0.68096892098918
$ws.Range("G5").Value = 0.002427531770031495
$ws.Range("I5").Value = 0.6234249009485318
$ws.Range("L5").Value = 0.3543873542648157
$ws.Range("M5").Value = 8.406049450756711
$ws.Range("O5").Value = 2.34957455767578
$ws.Range("C6").Value = 0.1450926412085494
$ws.Range("D6").Value = 0.0241467335481147
$ws.Range("E6").Value = 0.177827770774222
$ws.Range("F6").Value = 0.6817036818244446
$ws.Range("G6").Value = 0.002427767686527369
$ws.Range("I6").Value = 0.6246625006658064
$ws.Range("L6").Value = 0.3525423131965795
$ws.Range("M6").Value = 8.350422435628275
$ws.Range("O6").Value = 2.353038612139301
$ws.Range("C7").Value = 0.1448793551761582
$ws.Range("D7").Value = 0.02439305156950411
$ws.Range("E7").Value = 0.1824739362423813
$ws.Range("F7").Value = 0.676709214947536
$ws.Range("G7").Value = 0.002426142582229052
$ws.Range("I7").Value = 0.6161794299590682
$ws.Range("L7").Value = 0.3653708197547019
$ws.Range("M7").Value = 8.736242209425654
$ws.Range("O7").Value = 2.329375418769672
$ws.Range("C8").Value = 0.1443480659625891
$ws.Range("D8").Value = 0.02552308135666692
$ws.Range("E8").Value = 0.2033714891905731
$ws.Range("F8").Value = 0.6573173335915143
$ws.Range("G8").Value = 0.002419277254308798
$ws.Range("I8").Value = 0.5814352747116196
$ws.Range("L8").Value = 0.4225456257628935
$ws.Range("M8").Value = 10.43173342170832
$ws.Range("O8").Value = 2.234514714380879
$ws.Range("C9").Value = 0.1448319889512248
$ws.Range("D9").Value = 0.02790013214357856
$ws.Range("E9").Value = 0.2458906355087862
$ws.Range("F9").Value = 0.6291465358328949
$ws.Range("G9").Value = 0.002406900975727609
$ws.Range("I9").Value = 0.5233232911064896
$ws.Range("L9").Value = 0.5369496725034821
$ws.Range("M9").Value = 13.73362012078138
$ws.Range("O9").Value = 2.084099889397578
$ws.Range("C10").Value = 0.1461708331080018
$ws.Range("D10").Value = 0.02975160189743775
$ws.Range("E10").Value = 0.2781946997666438
$ws.Range("F10").Value = 0.614777789658298
$ws.Range("G10").Value = 0.002398457665094463
$ws.Range("I10").Value = 0.4870673951109197
$ws.Range("L10").Value = 0.6226673955724493
$ws.Range("M10").Value = 16.1484756729497
$ws.Range("O10").Value = 1.996463043686163
$ws.Range("C11").Value = 0.1470059100821999
$ws.Range("D11").Value = 0.03061838413434259
$ws.Range("E11").Value = 0.293152523338108
$ws.Range("F11").Value = 0.6096856131476329
$ws.Range("G11").Value = 0.002394754712894955
$ws.Range("I11").Value = 0.4720445987165043
$ws.Range("L11").Value = 0.6620901594043289
$ws.Range("M11").Value = 17.24542108111217
$ws.Range("O11").Value = 1.961818476011189
$ws.Range("C12").Value = 0.147355730266014
$ws.Range("D12").Value = 0.03095028952370171
$ws.Range("E12").Value = 0.2988570760527551
$ws.Range("F12").Value = 0.6079709569268203
$ws.Range("G12").Value = 0.002393372102165996
$ws.Range("I12").Value = 0.4665734772925347
$ws.Range("L12").Value = 0.6770860080498551
$ws.Range("M12").Value = 17.66064323457852
$ws.Range("O12").Value = 1.949472171871207
$ws.Range("C13").Value = 0.1472788780662313
$ws.Range("D13").Value = 0.03087864193729928
$ws.Range("E13").Value = 0.2976266608064293
$ws.Range("F13").Value = 0.6083306408207818
$ws.Range("G13").Value = 0.002393669003342279
$ws.Range("I13").Value = 0.4677419981020599
$ws.Range("L13").Value = 0.6738532976227134
$ws.Range("M13").Value = 17.57122412455919
$ws.Range("O13").Value = 1.952096442439569
$ws.Range("C14").Value = 0.1470340102453918
$ws.Range("D14").Value = 0.03064561569208735
$ws.Range("E14").Value = 0.293621016170917
$ws.Range("F14").Value = 0.6095402310829456
$ws.Range("G14").Value = 0.002394640572618861
$ws.Range("I14").Value = 0.4715900870766419
$ws.Range("L14").Value = 0.6633224947736664
$ws.Range("M14").Value = 17.27958464114329
$ws.Range("O14").Value = 1.960787121854679
$ws.Range("C15").Value = 0.1468884300817024
$ws.Range("D15").Value = 0.03050336331826031
$ws.Range("E15").Value = 0.291172777883034
$ws.Range("F15").Value = 0.6103091414208279
$ws.Range("G15").Value = 0.002395238235252661
$ws.Range("I15").Value = 0.4739756904313879
$ws.Range("L15").Value = 0.6568810106600154
$ws.Range("M15").Value = 17.10092715497512
$ws.Range("O15").Value = 1.966211720190586
$ws.Range("C16").Value = 0.1461208991286185
$ws.Range("D16").Value = 0.02969546157663672
$ws.Range("E16").Value = 0.2772226718562649
$ws.Range("F16").Value = 0.6151401438159922
$ws.Range("G16").Value = 0.002398702418894132
$ws.Range("I16").Value = 0.4880792872334361
$ws.Range("L16").Value = 0.6201001315759811
$ws.Range("M16").Value = 16.07676086843912
$ws.Range("O16").Value = 1.998834089955096
$ws.Range("C17").Value = 0.1457086934457976
$ws.Range("D17").Value = 0.02920622480462498
$ws.Range("E17").Value = 0.2687337701718633
$ws.Range("F17").Value = 0.6184779132337752
$ws.Range("G17").Value = 0.002400862761626149
$ws.Range("I17").Value = 0.4971123574687368
$ws.Range("L17").Value = 0.5976500560104796
$ws.Range("M17").Value = 15.44810340769357
$ws.Range("O17").Value = 2.020199767464618
$ws.Range("C18").Value = 0.1454928097394799
$ws.Range("D18").Value = 0.02892712915556928
$ws.Range("E18").Value = 0.2638757244501306
$ws.Range("F18").Value = 0.6205332847038321
$ws.Range("G18").Value = 0.002402118331870253
$ws.Range("I18").Value = 0.5024459385003865
$ws.Range("L18").Value = 0.5847774216163941
$ws.Range("M18").Value = 15.08636312990274
$ws.Range("O18").Value = 2.032978907565052
$ws.Range("C19").Value = 0.1454233251003672
$ws.Range("D19").Value = 0.02883302302328161
$ws.Range("E19").Value = 0.2622350131382944
$ws.Range("F19").Value = 0.6212523036731099
$ws.Range("G19").Value = 0.002402545685673555
$ws.Range("I19").Value = 0.50427530983983
$ws.Range("L19").Value = 0.5804256850630907
$ws.Range("M19").Value = 14.9638560385855
$ws.Range("O19").Value = 2.037389230683317
$ws.Range("C20").Value = 0.1457503709751364
$ws.Range("D20").Value = 0.02925806579641232
$ws.Range("E20").Value = 0.2696348695712771
$ws.Range("F20").Value = 0.6181085277319056
$ws.Range("G20").Value = 0.002400631445814328
$ws.Range("I20").Value = 0.4961364444750487
$ws.Range("L20").Value = 0.6000357202286182
$ws.Range("M20").Value = 15.51504049658035
$ws.Range("O20").Value = 2.017874466731143
$ws.Range("C21").Value = 0.1471050131061276
$ws.Range("D21").Value = 0.03071396023079842
$ws.Range("E21").Value = 0.2947964544697612
$ws.Range("F21").Value = 0.6091790988146357
$ws.Range("G21").Value = 0.002394354668665886
$ws.Range("I21").Value = 0.4704538514130086
$ws.Range("L21").Value = 0.666413773960727
$ws.Range("M21").Value = 17.36525025374425
$ws.Range("O21").Value = 1.958213306637816
$ws.Range("C22").Value = 0.1481867157846182
$ws.Range("D22").Value = 0.03168696039140428
$ws.Range("E22").Value = 0.3114775225005673
$ws.Range("F22").Value = 0.604591320887863
$ws.Range("G22").Value = 0.002390366666102232
$ws.Range("I22").Value = 0.45494064802417
$ws.Range("L22").Value = 0.7101908282012062
$ws.Range("M22").Value = 18.57352677411563
$ws.Range("O22").Value = 1.923736540071161
$ws.Range("C23").Value = 0.1475910466172508
$ws.Range("D23").Value = 0.03116563476696399
$ws.Range("E23").Value = 0.3025519779543089
$ws.Range("F23").Value = 0.6069236979416104
$ws.Range("G23").Value = 0.002392484760453666
$ws.Range("I23").Value = 0.4631018503257565
$ws.Range("L23").Value = 0.6867880929932824
$ws.Range("M23").Value = 17.92871150315284
$ws.Range("O23").Value = 1.941716915297235
$ws.Range("C24").Value = 0.1457314629279693
$ws.Range("D24").Value = 0.02923462172613256
$ws.Range("E24").Value = 0.2692274131525494
$ws.Range("F24").Value = 0.6182751023763373
$ws.Range("G24").Value = 0.002400735981407652
$ws.Range("I24").Value = 0.496577218058075
$ws.Range("L24").Value = 0.5989570552855525
$ws.Range("M24").Value = 15.48477922974138
$ws.Range("O24").Value = 2.018924192767258
$ws.Range("C25").Value = 0.1445324253274833
$ws.Range("D25").Value = 0.0272391612243581
$ws.Range("E25").Value = 0.2342107458516125
$ws.Range("F25").Value = 0.6356795813643217
$ws.Range("G25").Value = 0.002410134065591072
$ws.Range("I25").Value = 0.5379375901237005
$ws.Range("L25").Value = 0.5057288935074382
$ws.Range("M25").Value = 12.84254550745828
$ws.Range("O25").Value = 2.120859903334946
